$d = $word.ActiveDocument

# The bullet's final sentence currently lives in a single run:
#   " Helped to start the UW-FOSS club by helping with general logistics
#     (email and meetings) and working on the website with the founder."
# It needs to become five runs (same formatting throughout) with the
# wording "helping with" -> "facilitating" and "working on" -> "developing".

$full = $d.Content
$foundFull = $full.Find.Execute( `
    " Helped to start the UW-FOSS club by helping with general logistics (email and meetings) and working on the website with the founder.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($foundFull) {
    $base = $full.Start

    # Sub-ranges (offsets relative to $base) of the five target segments,
    # derived from the original sentence text.
    $seg1 = $d.Range($base+37,  $base+49)   # "helping with"  -> "facilitating"
    $seg2 = $d.Range($base+49,  $base+93)   # " general logistics (email and meetings) and "
    $seg3 = $d.Range($base+93,  $base+103)  # "working on"    -> "developing"
    $seg4 = $d.Range($base+103, $base+133)  # " the website with the founder."

    # Step 1: mark every segment with a (temporary) formatting toggle so
    # each becomes a genuinely distinct run from its neighbours before any
    # text is rewritten -- otherwise a text-only edit on one segment gets
    # re-coalesced into an adjacent, already-identical-looking run.
    $seg1.Font.Bold = $true
    $seg2.Font.Bold = $true
    $seg3.Font.Bold = $true
    $seg4.Font.Bold = $true

    # Step 2: rewrite the text of each segment (segments 2 and 4 keep their
    # original wording, but still need to be "touched" to end up as their
    # own runs).
    $seg1.Text = "facilitating"
    $seg2.Text = " general logistics (email and meetings) and "
    $seg3.Text = "developing"
    $seg4.Text = " the website with the founder."

    # Step 3: drop the temporary formatting again so the new runs end up
    # with formatting identical to the untouched first segment.
    $seg1.Font.Bold = $false
    $seg2.Font.Bold = $false
    $seg3.Font.Bold = $false
    $seg4.Font.Bold = $false
}
